# Insert a new data row at row 314 (pushing existing rows 314..365 down to
# 315..366) and populate the newly inserted row with the new observation.
#
# Resulting effects:
#   - Sheet dimension grows from A1:R365 to A1:R366
#   - Old row 314 data moves to row 315, old row 315 -> 316, ... old row 365 -> 366
#   - New row 314 receives a fresh record

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 314..365 down by inserting a new blank row at 314.
$ws.Rows.Item(314).Insert()

# Populate the newly inserted row 314 with the new record values.
$ws.Cells.Item(314, 1).Value = 7
$ws.Cells.Item(314, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(314, 3).Value = "Ñuble"
$ws.Cells.Item(314, 4).Value = 45218
$ws.Cells.Item(314, 5).Value = 16
$ws.Cells.Item(314, 6).Value = 100112032
$ws.Cells.Item(314, 7).Value = "Zapallo italiano"
$ws.Cells.Item(314, 8).Value = "Sin especificar"
$ws.Cells.Item(314, 9).Value = "Primera"
$ws.Cells.Item(314, 10).Value = 100
$ws.Cells.Item(314, 11).Value = 15000
$ws.Cells.Item(314, 12).Value = 15000
$ws.Cells.Item(314, 13).Value = 15000
$ws.Cells.Item(314, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(314, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(314, 16).Value = 300
$ws.Cells.Item(314, 17).Value = 50
$ws.Cells.Item(314, 18).Value = "Hortaliza"

# Make sure the date cell keeps the same date/time number format used by the
# rest of column D (style already carried over from Insert, but set it
# explicitly to be safe).
$ws.Cells.Item(314, 4).NumberFormat = $ws.Cells.Item(315, 4).NumberFormat
